# "updated BOM with R18 DNP"
#
# R18 (100 ohm resistor) is marked as Do-Not-Populate on the BOM sheet, and
# since it will not be populated, it no longer needs to be purchased, so its
# line item is removed from the "DK Order" (DigiKey order) sheet.

$wb = $excel.ActiveWorkbook
$wsBOM = $wb.Worksheets.Item("BOM")
$wsDK  = $wb.Worksheets.Item("DK Order")

# BOM sheet: R18 is row 34 -> flag it DNP in column F.
$wsBOM.Range("F34").Value = "DNP"

# DK Order sheet: R18's order line is row 33 -> remove it entirely, shifting
# every subsequent row up by one.
$wsDK.Rows.Item(33).Delete()

# Leave the selections/active sheet the way the author left the workbook:
# cursor resting on the (now shifted-up) row 33 in "DK Order", and "BOM" as
# the active/front sheet with the cursor near the bottom of its table.
$wsDK.Rows.Item(33).Select()
$wsBOM.Select()
$wsBOM.Range("I44").Select()
